$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to Text format so that
# numeric-looking strings (e.g. "1.003") are stored as text, matching
# the original inline-string cell contents, then strip the temporary
# number-format style back off so cell styling is unchanged.
$priceRange = $ws.Range("D2:D51")
$volRange = $ws.Range("E2:E51")
$priceRange.NumberFormat = "@"
$volRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.483.77"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "1.848.89"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "258.30"
$ws.Range("E5").Value = "  -7.24%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "0.5237"
$ws.Range("E7").Value = "  +2.73%  "
$ws.Range("D8").Value = "0.3291"
$ws.Range("E8").Value = "  -5.85%  "
$ws.Range("D9").Value = "0.06732"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "19.42"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").Value = "0.7786"
$ws.Range("E11").Value = "  -3.35%  "
$ws.Range("D12").Value = "0.07693"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "1.827.37"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "88.38"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "5.055"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "14.17"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "0.000007897"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("D20").Value = "26.520.31"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").Value = "2.099.05"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").Value = "4.613"
$ws.Range("D23").Value = "9.717"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("D25").Value = "2.378"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "144.13"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "1.650"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "17.04"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "111.55"
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "4.229"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("D31").Value = "4.203"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "0.08800"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "0.04880"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "1.144"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("D35").Value = "2.860"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "0.7089"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").Value = "3.092"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").Value = "0.01812"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "0.4965"
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("D40").Value = "2.214"
$ws.Range("E40").Value = "  -7.30%  "
$ws.Range("D41").Value = "114.39"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "0.9033"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("D43").Value = "6.080"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "7.812"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").Value = "0.4306"
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("E47").Value = "  -5.20%  "
$ws.Range("D48").Value = "9.204"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").Value = "0.05920"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "35.32"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.438"
$ws.Range("E51").Value = "  -3.39%  "

# Remove the temporary text-format styling so the cells keep their
# original (unstyled) appearance.
$priceRange.ClearFormats()
$volRange.ClearFormats()
